$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("material-diameter")

$ws.Range("B20").Value = "wwmd"
$ws.Range("C20").Value = -0.6625668846702162
$ws.Range("D20").Value = 0.001545301883302395
$ws.Range("E20").Value = 0.0009759820188660768
$ws.Range("F20").Value = 10724.66072463989
$ws.Range("G20").Value = -0.02545172565787205
$ws.Range("H20").Value = -0.03511352692945889
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = -0.07780546049260043
